# Slide 5, "Content Placeholder 2" (Shapes.Item(2)): the first bullet-less
# paragraph reads "First, download the zip file Lab4.zip". The zip file was
# renamed, so the sentence becomes "First, download the zip file
# Lab3-Fa2019.zip", typed in as three separate runs (matching the authored
# edit): "First, download the zip ", "file Lab3-Fa2019.", and "zip".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)
$para = $shape.TextFrame.TextRange.Paragraphs(1)

$para.Text = "First, download the zip "
$run2 = $para.InsertAfter("file Lab3-Fa2019.")
$run3 = $para.InsertAfter("zip")
